$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A11").Value = "Login with valid username and password"
$ws.Range("B11").Value = "PASSED"
$ws.Range("C11").Value = "chrome"

$ws.Range("A12").Value = "Login with valid username and password"
$ws.Range("B12").Value = "PASSED"
$ws.Range("C12").Value = "chrome"
